# Applies the "Keep alive" / "Statistics" rework described in the commit.
#
# Summary of changes:
#  1. Slide 11 ("Disconnect computer") -> "Keep alive": new title, new body
#     copy (2 paragraphs) and the screenshot picture is repositioned.
#  2. Slide 12 ("Statistics") -> "Statistics - Processing stage": title
#     gains a suffix, the empty content placeholder is replaced by a
#     2x4 results table plus a free-floating textbox with two summary
#     bullets.
#  3. A brand-new slide ("Statistics - Tree construct stage") is inserted
#     right after the former Statistics slide, using the same
#     "Title and Content" layout, title filled in, body left empty.

$p = $ppt.ActivePresentation

# EMU -> point helper (COM measurements are expressed in points).
function EmuToPt([double]$emu) { return $emu / 12700.0 }

# ---------------------------------------------------------------------
# 1. Slide 11: "Disconnect computer" -> "Keep alive"
# ---------------------------------------------------------------------
$slideKeepAlive = $p.Slides.Item(11)

$titleShape = $slideKeepAlive.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Keep alive"

$bodyShape = $slideKeepAlive.Shapes.Item(2)
$bodyShape.TextFrame.TextRange.Text = "The master constantly checks how many computers are running with keep-alive message." + "`r" + "If a computer crashes, the master change the responsibility of the computer that crashed to another computer."

$picShape = $slideKeepAlive.Shapes.Item(3)
$picShape.Left = EmuToPt 6219932
$picShape.Top = EmuToPt 3977220

# ---------------------------------------------------------------------
# 2. Slide 12: "Statistics" -> "Statistics - Processing stage"
# ---------------------------------------------------------------------
$slideStats = $p.Slides.Item(12)

$statsTitle = $slideStats.Shapes.Item(1)
$statsTitle.TextFrame.TextRange.Text = "Statistics - Processing stage"

# Remove the (empty) content placeholder; it gets replaced by a table and
# a separate textbox below.
$slideStats.Shapes.Item(2).Delete()

$tableShape = $slideStats.Shapes.AddTable(4, 2, (EmuToPt 2656114), (EmuToPt 2725054), (EmuToPt 5369560), (EmuToPt 1411516))
$tableShape.Name = "Content Placeholder 3"

$tbl = $tableShape.Table
$tbl.Columns.Item(1).Width = EmuToPt 2141791
$tbl.Columns.Item(2).Width = EmuToPt 3227769
$tbl.Rows.Item(1).Height = EmuToPt 352879
$tbl.Rows.Item(2).Height = EmuToPt 352879
$tbl.Rows.Item(3).Height = EmuToPt 352879
$tbl.Rows.Item(4).Height = EmuToPt 352879

$tableData = @(
    @(" ", "Average time after 20 test"),
    @("1 computer 1 worker", "4.7475 sec"),
    @("1 computer 4 workers", "4.8825 sec"),
    @("5 computers 4 workers", "4.79 sec")
)

for ($r = 1; $r -le 4; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cellRange = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
        $cellRange.Text = $tableData[$r - 1][$c - 1]
        $cellRange.Font.Size = 16
        $cellRange.ParagraphFormat.Alignment = 2
    }
}

$statsTextBox = $slideStats.Shapes.AddTextbox(1, (EmuToPt 1295401), (EmuToPt 4136570), (EmuToPt 9601196), (EmuToPt 1739298))
$statsTextBox.Name = "Content Placeholder 2"
$statsTextBox.TextFrame.TextRange.Text = "As the size of the big data increases, it will be better to process with a large number of computers." + "`r" + "As the size of the big data decreases, network messaging time will be more significant."
$statsTextBox.TextFrame.TextRange.ParagraphFormat.Alignment = 1

# ---------------------------------------------------------------------
# 3. Insert the new "Statistics - Tree construct stage" slide right
#    after the Statistics slide (now at position 12).
# ---------------------------------------------------------------------
$layout = $slideStats.CustomLayout
$treeSlide = $p.Slides.AddSlide(13, $layout)
$treeSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Statistics - Tree construct stage"
